# Refresh the "Price" (D) and "Volume(1h)" (E) columns on the cryptos sheet
# with the latest scraped figures. Values that look like a plain decimal
# number are prefixed with a leading apostrophe so Excel stores them as
# literal text (matching the original inlineStr cells) instead of silently
# re-typing them as numbers and dropping trailing zeros / precision.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.757.87'
$ws.Range('E2').Value = '  +0.39%  '
$ws.Range('D3').Value = '2.835.30'
$ws.Range('E3').Value = '  +2.55%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '''353.55'
$ws.Range('E5').Value = '  +5.90%  '
$ws.Range('D6').Value = '''113.53'
$ws.Range('E6').Value = '  -2.43%  '
$ws.Range('D7').Value = '''0.565'
$ws.Range('E7').Value = '  +4.99%  '
$ws.Range('D8').Value = '''1.00'
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').Value = '''0.600'
$ws.Range('E9').Value = '  +4.49%  '
$ws.Range('E10').Value = '  -0.40%  '
$ws.Range('D11').Value = '''0.0853'
$ws.Range('E11').Value = '  -1.73%  '
$ws.Range('D12').Value = '''20.07'
$ws.Range('E12').Value = '  -1.10%  '
$ws.Range('E13').Value = '  +1.10%  '
$ws.Range('D14').Value = '''7.72'
$ws.Range('E14').Value = '  +0.70%  '
$ws.Range('D15').Value = '3.284.82'
$ws.Range('E15').Value = '  +2.84%  '
$ws.Range('D16').Value = '2.830.68'
$ws.Range('E16').Value = '  +2.28%  '
$ws.Range('E17').Value = '  +0.46%  '
$ws.Range('D18').Value = '51.619.64'
$ws.Range('E18').Value = '  +0.11%  '
$ws.Range('D19').Value = '''7.38'
$ws.Range('E19').Value = '  +7.34%  '
$ws.Range('E20').Value = '  -3.24%  '
$ws.Range('E21').Value = '  -0.01%  '
$ws.Range('D22').Value = '0.0₃0995'
$ws.Range('E22').Value = '  +1.94%  '
$ws.Range('D23').Value = '''270.91'
$ws.Range('E23').Value = '  -2.75%  '
$ws.Range('D24').Value = '''69.76'
$ws.Range('E24').Value = '  +0.08%  '
$ws.Range('D25').Value = '''2.78'
$ws.Range('E25').Value = '  +3.64%  '
$ws.Range('E26').Value = '  -0.30%  '
$ws.Range('E27').Value = '  +0.14%  '
$ws.Range('D28').Value = '''10.32'
$ws.Range('E28').Value = '  +1.39%  '
$ws.Range('E29').Value = '  +1.12%  '
$ws.Range('D30').Value = '''0.139'
$ws.Range('E30').Value = '  -1.51%  '
$ws.Range('D31').Value = '''34.15'
$ws.Range('E31').Value = '  -2.67%  '
$ws.Range('D32').Value = '''50.67'
$ws.Range('E32').Value = '  +1.21%  '
$ws.Range('D33').Value = '''5.83'
$ws.Range('E33').Value = '  +4.63%  '
$ws.Range('E34').Value = '  +25.55%  '
$ws.Range('D35').Value = '''0.0826'
$ws.Range('E35').Value = '  +0.26%  '
$ws.Range('D36').Value = '''0.999'
$ws.Range('E36').Value = '  -0.19%  '
$ws.Range('E37').Value = '  -0.34%  '
$ws.Range('D38').Value = '''4.89'
$ws.Range('E38').Value = '  -2.20%  '
$ws.Range('D39').Value = '''3.20'
$ws.Range('E39').Value = '  -1.31%  '
$ws.Range('D40').Value = '''18.12'
$ws.Range('E40').Value = '  -4.66%  '
$ws.Range('D41').Value = '''23.60'
$ws.Range('E41').Value = '  +1.53%  '
$ws.Range('E42').Value = '  +1.98%  '
$ws.Range('D43').Value = '''2.54'
$ws.Range('E43').Value = '  +3.87%  '
$ws.Range('D44').Value = '''125.71'
$ws.Range('E44').Value = '  -1.07%  '
$ws.Range('E45').Value = '  -0.25%  '
$ws.Range('D46').Value = '2.085.06'
$ws.Range('E46').Value = '  -0.29%  '
$ws.Range('E47').Value = '  +0.30%  '
$ws.Range('E48').Value = '  +3.66%  '
$ws.Range('E49').Value = '  +2.95%  '
$ws.Range('E50').Value = '  +6.68%  '
$ws.Range('D51').Value = '''60.93'
$ws.Range('E51').Value = '  +1.32%  '
